$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 975.7143
$ws.Range("J51").Value = 988.3333
$ws.Range("L51").Value = 988.3333
$ws.Range("N51").Value = -1956.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 31551.428
$ws.Range("J128").Value = 31551.428
$ws.Range("L128").Value = 31551.428
$ws.Range("N128").Value = -41511.428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 39401.668
$ws.Range("J130").Value = 39401.668
$ws.Range("L130").Value = 39401.668
$ws.Range("N130").Value = -49441.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4055468.8
$ws.Range("I137").Value = 1924440.1
$ws.Range("J137").Value = 9092446
$ws.Range("K137").Value = 5773320.300000001
$ws.Range("L137").Value = 27277338
$ws.Range("M137").Value = -5770770.300000001
$ws.Range("N137").Value = -27282438

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 297957.22
$ws.Range("I138").Value = 1131.2941
$ws.Range("J138").Value = 1215419.1
$ws.Range("K138").Value = 3393.8823
$ws.Range("L138").Value = 3646257.3
$ws.Range("M138").Value = 1746.1177
$ws.Range("N138").Value = -3656537.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 10032.857
$ws.Range("J21").Value = 15333.333
$ws.Range("L21").Value = 15333.333
$ws.Range("N21").Value = -16081.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2984.16
$ws.Range("I32").Value = 2840.9795
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 2840.9795
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -2553.9795
$ws.Range("N32").Value = -10574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 13364.28
$ws.Range("J37").Value = 13459.409
$ws.Range("L37").Value = 13459.409
$ws.Range("N37").Value = -14005.409

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3580.16
$ws.Range("I45").Value = 2771
$ws.Range("J45").Value = 4456.75
$ws.Range("K45").Value = 2771
$ws.Range("L45").Value = 4456.75
$ws.Range("M45").Value = -2394
$ws.Range("N45").Value = -5210.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1312.7333
$ws.Range("I61").Value = 1306.6586
$ws.Range("J61").Value = 1375
$ws.Range("K61").Value = 1306.6586
$ws.Range("L61").Value = 1375
$ws.Range("M61").Value = -1094.6586
$ws.Range("N61").Value = -1799

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 20490.666
$ws.Range("J80").Value = 27236
$ws.Range("L80").Value = 27236
$ws.Range("N80").Value = -29232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 20490.666
$ws.Range("J83").Value = 27236
$ws.Range("L83").Value = 81708
$ws.Range("N83").Value = -91692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 13066.667
$ws.Range("J109").Value = 13066.667
$ws.Range("L109").Value = 13066.667
$ws.Range("N109").Value = -15840.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1943.862
$ws.Range("I122").Value = 1810.0741
$ws.Range("K122").Value = 5430.2223
$ws.Range("M122").Value = -2980.2223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 53337.5
$ws.Range("J128").Value = 53337.5
$ws.Range("L128").Value = 53337.5
$ws.Range("N128").Value = -63297.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 48784.625
$ws.Range("J129").Value = 48784.625
$ws.Range("L129").Value = 48784.625
$ws.Range("N129").Value = -58784.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1312.7333
$ws.Range("I136").Value = 1306.6586
$ws.Range("J136").Value = 1375
$ws.Range("K136").Value = 3919.9758
$ws.Range("L136").Value = 4125
$ws.Range("M136").Value = -1369.9758
$ws.Range("N136").Value = -9225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 47049.727
$ws.Range("I82").Value = 82878.5
$ws.Range("K82").Value = 82878.5
$ws.Range("M82").Value = -82495.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 47049.727
$ws.Range("I85").Value = 82878.5
$ws.Range("K85").Value = 82878.5
$ws.Range("M85").Value = -81552.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47799.4
$ws.Range("J20").Value = 47799.4
$ws.Range("L20").Value = 47799.4
$ws.Range("N20").Value = -48271.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 47799.4
$ws.Range("J30").Value = 47799.4
$ws.Range("L30").Value = 47799.4
$ws.Range("N30").Value = -47981.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2110.2432
$ws.Range("I31").Value = 1668.0869
$ws.Range("J31").Value = 2836.6428
$ws.Range("K31").Value = 1668.0869
$ws.Range("L31").Value = 2836.6428
$ws.Range("M31").Value = -1373.0869
$ws.Range("N31").Value = -3426.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2110.2432
$ws.Range("I34").Value = 1668.0869
$ws.Range("J34").Value = 2836.6428
$ws.Range("K34").Value = 1668.0869
$ws.Range("L34").Value = 2836.6428
$ws.Range("M34").Value = -1466.0869
$ws.Range("N34").Value = -3240.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4491.4287
$ws.Range("I105").Value = 4355.6665
$ws.Range("J105").Value = 4735.8
$ws.Range("K105").Value = 4355.6665
$ws.Range("L105").Value = 4735.8
$ws.Range("M105").Value = -2608.6665
$ws.Range("N105").Value = -8229.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 10780
$ws.Range("J109").Value = 10866.667
$ws.Range("L109").Value = 10866.667
$ws.Range("N109").Value = -12946.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 47799.4
$ws.Range("J128").Value = 47799.4
$ws.Range("L128").Value = 47799.4
$ws.Range("N128").Value = -57759.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 35326
$ws.Range("J131").Value = 35326
$ws.Range("L131").Value = 35326
$ws.Range("N131").Value = -45406

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3117.0417
$ws.Range("I132").Value = 2284.1052
$ws.Range("J132").Value = 6282.2
$ws.Range("K132").Value = 6852.3156
$ws.Range("L132").Value = 18846.6
$ws.Range("M132").Value = -4322.3156
$ws.Range("N132").Value = -23906.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 39450
$ws.Range("J133").Value = 39450
$ws.Range("L133").Value = 39450
$ws.Range("N133").Value = -44510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 49242.855
$ws.Range("J135").Value = 49242.855
$ws.Range("L135").Value = 49242.855
$ws.Range("N135").Value = -59382.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 34621.25
$ws.Range("J62").Value = 34621.25
$ws.Range("L62").Value = 34621.25
$ws.Range("N62").Value = -35993.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 34621.25
$ws.Range("J65").Value = 34621.25
$ws.Range("L65").Value = 103863.75
$ws.Range("N65").Value = -110727.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2137.2144
$ws.Range("I122").Value = 2131.5715
$ws.Range("J122").Value = 2142.8572
$ws.Range("K122").Value = 6394.7145
$ws.Range("L122").Value = 6428.571599999999
$ws.Range("M122").Value = -3944.7145
$ws.Range("N122").Value = -11328.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 27710.834
$ws.Range("J63").Value = 27710.834
$ws.Range("L63").Value = 27710.834
$ws.Range("N63").Value = -29208.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 27710.834
$ws.Range("J66").Value = 27710.834
$ws.Range("L66").Value = 83132.50199999999
$ws.Range("N66").Value = -90620.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3927.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3927.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11782.5
$ws.Range("N122").Value = -16682.5
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 23055.666
$ws.Range("J109").Value = 23055.666
$ws.Range("L109").Value = 23055.666
$ws.Range("N109").Value = -25829.666
